$d = $word.ActiveDocument

# Locate the field whose instruction text references the "userdoc" tag
# (e.g. ` m:userdoc 'zone1' `). This is the field that the
# TokenIteratorFieldRewriterSplit-based parser now emits as plain text
# runs instead of a real Word field.
$targetField = $null
foreach ($f in $d.Fields) {
    if ($f.Code.Text -match "userdoc") {
        $targetField = $f
    }
}

if ($targetField -ne $null) {
    # Find the paragraph that contains this field so we can replace the
    # whole paragraph content (fldChar begin / instrText* / fldChar end)
    # with two plain <w:t> runs holding the same text: "m" and
    # ":userdoc 'zone1'".
    $targetParagraph = $null
    foreach ($p in $d.Paragraphs) {
        if ($targetField.Code.Start -ge $p.Range.Start -and $targetField.Code.Start -lt $p.Range.End) {
            $targetParagraph = $p
        }
    }

    if ($targetParagraph -ne $null) {
        $rng = $targetParagraph.Range
        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:r><w:t>m</w:t></w:r>' +
               '<w:r><w:t>:userdoc ''zone1''</w:t></w:r>' +
               '</w:p>'
        $rng.InsertXML($xml)
    }
}
